$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "53.86") need to be
# forced to Text format first, otherwise Excel will auto-convert them into
# numeric values instead of keeping the original "price-as-text" formatting
# used throughout this sheet.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D20", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin data (prices / 1h volume changes, plus the two rows
# whose coins swapped rank position).
$ws.Range("D2").Value = "30.512.63"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.108.85"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "333.94"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.5266"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").Value = "0.4569"
$ws.Range("E8").Value = "  +4.50%  "
$ws.Range("D9").Value = "53.86"
$ws.Range("E9").Value = "  +14.53%  "
$ws.Range("D10").Value = "0.09023"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").Value = "1.186"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").Value = "24.52"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "2.102.16"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "6.806"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "7.859"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "96.93"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "0.00001132"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").Value = "19.50"
$ws.Range("E20").Value = "  +2.51%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "6.331"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "30.556.62"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").Value = "12.37"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").Value = "2.364"
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("D26").Value = "2.342.88"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").Value = "22.43"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").Value = "2.593"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "163.69"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "132.97"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "1.205"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("D32").Value = "0.1076"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").Value = "1.673"
$ws.Range("E33").Value = "  +8.28%  "
$ws.Range("D34").Value = "6.163"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "3.928"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("D36").Value = "10.48"
$ws.Range("E36").Value = "  +9.76%  "
$ws.Range("D37").Value = "0.02582"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "5.607"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.06842"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").Value = "12.81"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").Value = "0.2297"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").Value = "0.6931"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").Value = "1.245"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "2.364"
$ws.Range("E44").Value = "  +6.33%  "
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6405"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "14.00"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").Value = "3.659"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "0.00000000356"
$ws.Range("E49").Value = "  +28.44%  "
$ws.Range("D50").Value = "1.252"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").Value = "1.216"
$ws.Range("E51").Value = "  +0.74%  "
